# EMS_DB Design.xlsx edit:
# - Clear the "PendingUser(Data)" column (D27:D50) from the User(Data) table,
#   since that table is being removed from the design sheet.
# - Rename the "RoleUID" field (A32) to "RoleID" to match the actual DB column table.
# - Add four new User(Data) fields used for auth/password-reset support:
#   PasswordHash, HashKey, ResetPasswordToken, ResetPasswordTokenKey (rows 61-64).
# - Minor view/formatting touch-ups (row heights on the section-header rows,
#   and the active selection) left by the editing session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Clear the old "PendingUser(Data)" mirror column (D27:D50), keeping the
#    existing cell formatting/styles in place.
$ws.Range("D27:D50").ClearContents()

# 2. The old "RoleUID" row now reads "RoleID" (matches the DesignationID/RoleID
#    naming used in the detailed table further down the sheet).
$ws.Range("A32").Value = "RoleID"

# 3. Append the four new User(Data) fields below the existing list (row 60 was
#    the last populated row in that column).
$ws.Range("A61").Value = "PasswordHash"
$ws.Range("A62").Value = "HashKey"
$ws.Range("A63").Value = "ResetPasswordToken"
$ws.Range("A64").Value = "ResetPasswordTokenKey"

# Match the formatting of the row directly above (same column, same list).
$ws.Range("A60").Copy()
$ws.Range("A61:A64").PasteSpecial(-4122)
$ws.Range("A61:A64").RowHeight = 15.75
$excel.CutCopyMode = 0

# 4. The three section-header rows got a touch-up in row height.
$ws.Range("A2").EntireRow.RowHeight = 13
$ws.Range("A7").EntireRow.RowHeight = 13
$ws.Range("A10").EntireRow.RowHeight = 13

# 5. Leave the selection where the edits were made.
$ws.Range("A32").Select()
